$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row label updates
$ws.Range("A2").Value = "Julian Day of Year"
$ws.Range("A4").Value = "Unique Squirrels"

# Sleep per day (hrs)
$ws.Range("B6").Value = "13.37 ± 1.89"
$ws.Range("C6").Value = "12.61 ± 2.53"
$ws.Range("D6").Value = "13.30 ± 3.11"
$ws.Range("E6").Value = "10.77 ± 3.13"
$ws.Range("F6").Value = "12.66 ± 2.86"

# Sleep in daylight (hrs)
$ws.Range("B7").Value = "2.31 ± 0.86"
$ws.Range("C7").Value = "3.93 ± 1.92"
$ws.Range("D7").Value = "8.64 ± 2.82"
$ws.Range("E7").Value = "3.19 ± 2.60"
$ws.Range("F7").Value = "4.55 ± 3.36"

# Sleep in darkness (hrs)
$ws.Range("B8").Value = "11.07 ± 1.53"
$ws.Range("C8").Value = "8.69 ± 2.07"
$ws.Range("D8").Value = "4.66 ± 0.81"
$ws.Range("E8").Value = "7.58 ± 1.30"
$ws.Range("F8").Value = "8.10 ± 2.88"

# Total sleep transitions
$ws.Range("B9").Value = "453 ± 87"
$ws.Range("C9").Value = "363 ± 94"
$ws.Range("D9").Value = "313 ± 56"
$ws.Range("E9").Value = "283 ± 74"
$ws.Range("F9").Value = "361 ± 104"

# Sleep transitions in daylight
$ws.Range("B10").Value = "94 ± 27"
$ws.Range("C10").Value = "139 ± 59"
$ws.Range("D10").Value = "244 ± 63"
$ws.Range("E10").Value = "122 ± 83"
$ws.Range("F10").Value = "150 ± 84"

# Sleep transitions in darkness
$ws.Range("B11").Value = "359 ± 78`n"
$ws.Range("C11").Value = "223 ± 91`n"
$ws.Range("D11").Value = "69 ± 39`n"
$ws.Range("E11").Value = "161 ± 68`n"
$ws.Range("F11").Value = "211 ± 135"
